$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Content.Find.Execute("2023-09-01 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-09-02 Saturday", 2)

# Update the table of division problems directly by cell (row, col) to avoid
# ambiguity from duplicate values that appear more than once in the table.
$table = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="70÷4=17, 2"},
    @{Row=1;  Col=2; Text="82÷2=41, 0"},
    @{Row=1;  Col=3; Text="72÷2=36, 0"},
    @{Row=1;  Col=4; Text="79÷4=19, 3"},
    @{Row=1;  Col=5; Text="27÷6=4, 3"},

    @{Row=5;  Col=1; Text="81÷4=20, 1"},
    @{Row=5;  Col=2; Text="45÷4=11, 1"},
    @{Row=5;  Col=3; Text="64÷9=7, 1"},
    @{Row=5;  Col=4; Text="27÷4=6, 3"},
    @{Row=5;  Col=5; Text="43÷3=14, 1"},

    @{Row=9;  Col=1; Text="44÷7=6, 2"},
    @{Row=9;  Col=2; Text="68÷5=13, 3"},
    @{Row=9;  Col=3; Text="28÷5=5, 3"},
    @{Row=9;  Col=4; Text="79÷8=9, 7"},
    @{Row=9;  Col=5; Text="94÷3=31, 1"},

    @{Row=13; Col=1; Text="14÷7=2, 0"},
    @{Row=13; Col=2; Text="69÷2=34, 1"},
    @{Row=13; Col=3; Text="57÷4=14, 1"},
    @{Row=13; Col=4; Text="84÷3=28, 0"},
    @{Row=13; Col=5; Text="57÷2=28, 1"},

    @{Row=17; Col=1; Text="77÷8=9, 5"},
    @{Row=17; Col=2; Text="45÷5=9, 0"},
    @{Row=17; Col=3; Text="45÷6=7, 3"},
    @{Row=17; Col=4; Text="46÷5=9, 1"},
    @{Row=17; Col=5; Text="81÷4=20, 1"}
)

foreach ($u in $updates) {
    $cell = $table.Cell($u.Row, $u.Col)
    $cellRange = $cell.Range
    $cellRange.MoveEnd(1, -1)
    $cellRange.Text = $u.Text
}
